$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01321787834167481
$ws.Range("C2").Value = 0.02446446418762207
$ws.Range("D2").Value = 0.005193376541137695
$ws.Range("E2").Value = 0.01822609901428223
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.0632962703704834
$ws.Range("H2").Value = 0.01850481033325195
$ws.Range("I2").Value = 0.01909732818603516
$ws.Range("J2").Value = 0.01494874954223633
$ws.Range("K2").Value = 0.03093433380126953
$ws.Range("L2").Value = 0.005259943008422851
$ws.Range("M2").Value = 0.01699624061584473
$ws.Range("B3").Value = 0.08840374946594239
$ws.Range("C3").Value = 0.0273953914642334
$ws.Range("D3").Value = 0.0202277660369873
$ws.Range("E3").Value = 0.01152606010437012
$ws.Range("F3").Value = 0.01385688781738281
$ws.Range("G3").Value = 0.004719305038452149
$ws.Range("H3").Value = 0.09096217155456543
$ws.Range("I3").Value = 0.02821483612060547
$ws.Range("J3").Value = 0.03766140937805176
$ws.Range("K3").Value = 0.01563277244567871
$ws.Range("L3").Value = 0.01390180587768555
$ws.Range("M3").Value = 0.009477376937866211
$ws.Range("B4").Value = 0.02509331703186035
$ws.Range("C4").Value = 0.01291303634643555
$ws.Range("D4").Value = 0.005937910079956055
$ws.Range("E4").Value = 0.006662607192993164
$ws.Range("F4").Value = 0.02828989028930664
$ws.Range("G4").Value = 0.006249809265136718
$ws.Range("H4").Value = 0.01590471267700195
$ws.Range("I4").Value = 0.01238737106323242
$ws.Range("J4").Value = 0.01097283363342285
$ws.Range("K4").Value = 0.006257772445678711
$ws.Range("L4").Value = 0.03135590553283692
$ws.Range("M4").Value = 0.003134727478027344
$ws.Range("B5").Value = 0.01251769065856934
$ws.Range("C5").Value = 0.01331238746643066
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.015704345703125
$ws.Range("H5").Value = 0.01048812866210938
$ws.Range("I5").Value = 0.01461496353149414
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.01531906127929688
$ws.Range("B6").Value = 0.3644277572631836
$ws.Range("C6").Value = 0.0188596248626709
$ws.Range("D6").Value = 0.329889440536499
$ws.Range("E6").Value = 0.02056207656860352
$ws.Range("F6").Value = 0.3833872795104981
$ws.Range("G6").Value = 0.01256566047668457
$ws.Range("H6").Value = 0.1308046340942383
$ws.Range("I6").Value = 0.01892623901367187
$ws.Range("J6").Value = 0.1223444938659668
$ws.Range("K6").Value = 0.01421184539794922
$ws.Range("L6").Value = 0.3930227756500244
$ws.Range("M6").Value = 0.01550831794738769
$ws.Range("B7").Value = 0.6131297588348389
$ws.Range("C7").Value = 0.06164984703063965
$ws.Range("D7").Value = 0.2285384654998779
$ws.Range("E7").Value = 0.02266936302185059
$ws.Range("F7").Value = 0.5068877696990967
$ws.Range("G7").Value = 0.02262983322143555
$ws.Range("H7").Value = 0.6459254741668701
$ws.Range("I7").Value = 0.06580033302307128
$ws.Range("J7").Value = 0.3862461090087891
$ws.Range("K7").Value = 0.04363632202148438
$ws.Range("L7").Value = 0.6442729949951171
$ws.Range("M7").Value = 0.02774744033813476
